$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 09:52"

# Row 11 - Rusia
$ws.Range("B11").Value = 106498
$ws.Range("C11").Value = 7099
$ws.Range("D11").Value = 11619
$ws.Range("E11").Value = 93806
$ws.Range("G11").Value = 101
$ws.Range("H11").Value = 1073

# Row 45 - Chequia
$ws.Range("F45").Value = 68

# Row 61 - Kazajistan
$ws.Range("B61").Value = 3273
$ws.Range("C61").Value = 135
$ws.Range("E61").Value = 2429

# Row 66 - Oman (name unchanged)
$ws.Range("B66").Value = 2348
$ws.Range("C66").Value = 74
$ws.Range("D66").Value = 495
$ws.Range("E66").Value = 1843

# Row 67 - now Afganistan (was Armenia) - inserted above Armenia
$ws.Range("A67").Value = "Afganistan"
$ws.Range("B67").Value = 2171
$ws.Range("C67").Value = 232
$ws.Range("D67").Value = 260
$ws.Range("E67").Value = 1847
$ws.Range("F67").Value = 7
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 64

# Row 68 - now Armenia (was Croacia)
$ws.Range("A68").Value = "Armenia"
$ws.Range("B68").Value = 2066
$ws.Range("C68").Value = 134
$ws.Range("D68").Value = 929
$ws.Range("E68").Value = 1105
$ws.Range("F68").Value = 10
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 32

# Row 69 - now Croacia (was Uzbekistan)
$ws.Range("A69").Value = "Croacia"
$ws.Range("B69").Value = 2062
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 1288
$ws.Range("E69").Value = 707
$ws.Range("F69").Value = 19
$ws.Range("H69").Value = 67

# Row 70 - now Uzbekistan (was Irak)
$ws.Range("A70").Value = "Uzbekistan"
$ws.Range("B70").Value = 2017
$ws.Range("C70").Value = 15
$ws.Range("D70").Value = 1096
$ws.Range("E70").Value = 912
$ws.Range("F70").Value = 8
$ws.Range("H70").Value = 9

# Row 71 - now Irak (was Afganistan, the original row that gets removed from its old slot)
$ws.Range("A71").Value = "Irak"
$ws.Range("B71").Value = 2003
$ws.Range("D71").Value = 1346
$ws.Range("E71").Value = 565
$ws.Range("F71").Value = 0
$ws.Range("H71").Value = 92

# Row 84 - Eslovaquia
$ws.Range("B84").Value = 1396
$ws.Range("C84").Value = 5
$ws.Range("D84").Value = 524
$ws.Range("E84").Value = 849
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 23

# Row 85 - Lituania
$ws.Range("B85").Value = 1385
$ws.Range("C85").Value = 10
$ws.Range("D85").Value = 589
$ws.Range("E85").Value = 751

# Row 103 - Sri Lanka
$ws.Range("D103").Value = 139
$ws.Range("E103").Value = 503

# Row 110
$ws.Range("B110").Value = 500
$ws.Range("C110").Value = 9
$ws.Range("D110").Value = 65
$ws.Range("E110").Value = 404
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 31

# Row 128 - Maldivas
$ws.Range("B128").Value = 280
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 262
